# Updates gh-pages "want-to-go" / price data generated at 456a3b4.
# Sheets: 展览 (exhibitions), 演出 (performances), 本地生活 (local life),
# 全部类型 (all types / combined).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 6077
$ws.Range("F10").Value = 709
$ws.Range("F11").Value = 2381
$ws.Range("F12").Value = 2381
$ws.Range("F13").Value = 11
$ws.Range("F14").Value = 1701
$ws.Range("F15").Value = 591
$ws.Range("F16").Value = 243
$ws.Range("F17").Value = 668
$ws.Range("F18").Value = 4850
$ws.Range("F19").Value = 136
$ws.Range("F20").Value = 56
$ws.Range("F21").Value = 684
$ws.Range("F23").Value = 844
$ws.Range("F26").Value = 30
$ws.Range("F27").Value = 2387
$ws.Range("F32").Value = 472
$ws.Range("F34").Value = 795
$ws.Range("F35").Value = 42
$ws.Range("F36").Value = 15
$ws.Range("F38").Value = 1359
$ws.Range("F39").Value = 1332

$ws = $wb.Worksheets.Item("演出")
$ws.Range("G4").Value = "不可售"
$ws.Range("G5").Value = "不可售"
$ws.Range("F14").Value = 116
$ws.Range("F16").Value = 76
$ws.Range("F20").Value = 321
$ws.Range("F22").Value = 511

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 818
$ws.Range("F4").Value = 228

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 818
$ws.Range("F6").Value = 228
$ws.Range("F7").Value = 6077
$ws.Range("F8").Value = 6077
$ws.Range("B10").Value = "'2024-10-25"
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = "北京·伦敦西区音乐剧明星演唱会-经典版"
$ws.Range("D10").Value = "西直门外大街135号（北京展览馆内） 北京展览馆剧场"
$ws.Range("E10").Value = "2024.10.25 19:30-10.26 21:30"
$ws.Range("F10").Value = 18
$ws.Range("G10").Value = 144
$ws.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=89359"
$ws.Range("I10").Value = "//i0.hdslb.com/bfs/openplatform/202407/PzPiEKUI1721114840552.jpeg"
$ws.Range("C11").Value = "北京·魔术脱口秀《三块巧克力》——“亦苦亦甜”快乐人生魔法SHOW"
$ws.Range("D11").Value = "三里屯SOHO下沉广场最南端,6号商场B1层 爱乐汇艺术空间"
$ws.Range("E11").Value = "2024.10.25 19:30-11.02 20:45"
$ws.Range("F11").Value = 2
$ws.Range("G11").Value = 126
$ws.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=93292"
$ws.Range("I11").Value = "//i1.hdslb.com/bfs/openplatform/202410/i0w53n9w1728620588393.jpeg"
$ws.Range("F20").Value = 2381
$ws.Range("F22").Value = 11
$ws.Range("F23").Value = 1701
$ws.Range("F24").Value = 116
$ws.Range("F25").Value = 591
$ws.Range("F26").Value = 243
$ws.Range("F27").Value = 668
$ws.Range("F28").Value = 4850
$ws.Range("F29").Value = 56
$ws.Range("F30").Value = 684
$ws.Range("F35").Value = 30
$ws.Range("F36").Value = 2387
$ws.Range("F39").Value = 472
$ws.Range("F43").Value = 511
$ws.Range("F44").Value = 795
$ws.Range("F45").Value = 42
$ws.Range("F46").Value = 15
$ws.Range("F48").Value = 1359
$ws.Range("F50").Value = 1332
